$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing space from the header "Nr. " -> "Nr."
$ws.Range("A1").Value = "Nr."

# Row heights recalculated (e.g. by opening/re-saving in a newer Excel which
# rewraps the long "Fachrichtung" text in column C at slightly different
# widths). Apply the resulting explicit row heights for the rows that grew.
$ws.Rows.Item(1).RowHeight = 33
$ws.Rows.Item(5).RowHeight = 46.5
$ws.Rows.Item(8).RowHeight = 33
$ws.Rows.Item(11).RowHeight = 33
$ws.Rows.Item(15).RowHeight = 33
$ws.Rows.Item(18).RowHeight = 33

$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 19.5
$ws.Rows.Item(4).RowHeight = 19.5
$ws.Rows.Item(6).RowHeight = 19.5
$ws.Rows.Item(7).RowHeight = 19.5
$ws.Rows.Item(9).RowHeight = 19.5
$ws.Rows.Item(10).RowHeight = 19.5
$ws.Rows.Item(12).RowHeight = 19.5
$ws.Rows.Item(13).RowHeight = 19.5
$ws.Rows.Item(14).RowHeight = 19.5
$ws.Rows.Item(16).RowHeight = 19.5
$ws.Rows.Item(17).RowHeight = 19.5

# Move the active selection to A2, as it is after the edit was saved.
[void]$ws.Range("A2").Select()
